$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 10

# A10 is a genuine (but empty) text cell in the source data, like A2:A9.
# A leading apostrophe forces Excel's "text" quote-prefix so the cell is
# actually created instead of being cleared by an empty assignment.
$ws.Cells.Item($row, 1).Value = "'"

$ws.Cells.Item($row, 2).Value = "محمود"

# C10 looks numeric ("222") but must stay text, matching the rest of the
# column (all existing quantity values are stored as text). The leading
# apostrophe forces a text interpretation without attaching a "Text"
# number-format override.
$ws.Cells.Item($row, 3).Value = "'222"

$ws.Cells.Item($row, 4).Value = "النصر"
$ws.Cells.Item($row, 5).Value = "الرحلة 3"
$ws.Cells.Item($row, 6).Value = "C2"
$ws.Cells.Item($row, 7).Value = "ABC"
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٧:٤١:٣٥ م"
